$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-03 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-04 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("815÷9=90, 5", $true, $false, $false, $false, $false, $true, 1, $false, "325÷5=65, 0", 2) | Out-Null
$d.Content.Find.Execute("763÷2=381, 1", $true, $false, $false, $false, $false, $true, 1, $false, "789÷3=263, 0", 2) | Out-Null
$d.Content.Find.Execute("305÷5=61, 0", $true, $false, $false, $false, $false, $true, 1, $false, "952÷2=476, 0", 2) | Out-Null
$d.Content.Find.Execute("122÷3=40, 2", $true, $false, $false, $false, $false, $true, 1, $false, "532÷2=266, 0", 2) | Out-Null
$d.Content.Find.Execute("705÷2=352, 1", $true, $false, $false, $false, $false, $true, 1, $false, "642÷7=91, 5", 2) | Out-Null
$d.Content.Find.Execute("651÷3=217, 0", $true, $false, $false, $false, $false, $true, 1, $false, "133÷3=44, 1", 2) | Out-Null
$d.Content.Find.Execute("937÷2=468, 1", $true, $false, $false, $false, $false, $true, 1, $false, "116÷2=58, 0", 2) | Out-Null
$d.Content.Find.Execute("396÷3=132, 0", $true, $false, $false, $false, $false, $true, 1, $false, "655÷2=327, 1", 2) | Out-Null
$d.Content.Find.Execute("202÷3=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "187÷4=46, 3", 2) | Out-Null
$d.Content.Find.Execute("691÷9=76, 7", $true, $false, $false, $false, $false, $true, 1, $false, "834÷8=104, 2", 2) | Out-Null
$d.Content.Find.Execute("960÷2=480, 0", $true, $false, $false, $false, $false, $true, 1, $false, "943÷9=104, 7", 2) | Out-Null
$d.Content.Find.Execute("294÷3=98, 0", $true, $false, $false, $false, $false, $true, 1, $false, "189÷7=27, 0", 2) | Out-Null
$d.Content.Find.Execute("453÷5=90, 3", $true, $false, $false, $false, $false, $true, 1, $false, "884÷8=110, 4", 2) | Out-Null
$d.Content.Find.Execute("155÷7=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "973÷7=139, 0", 2) | Out-Null
$d.Content.Find.Execute("512÷9=56, 8", $true, $false, $false, $false, $false, $true, 1, $false, "186÷2=93, 0", 2) | Out-Null
$d.Content.Find.Execute("669÷4=167, 1", $true, $false, $false, $false, $false, $true, 1, $false, "662÷5=132, 2", 2) | Out-Null
$d.Content.Find.Execute("394÷9=43, 7", $true, $false, $false, $false, $false, $true, 1, $false, "236÷3=78, 2", 2) | Out-Null
$d.Content.Find.Execute("402÷5=80, 2", $true, $false, $false, $false, $false, $true, 1, $false, "231÷5=46, 1", 2) | Out-Null
$d.Content.Find.Execute("159÷5=31, 4", $true, $false, $false, $false, $false, $true, 1, $false, "605÷4=151, 1", 2) | Out-Null
$d.Content.Find.Execute("113÷7=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "385÷4=96, 1", 2) | Out-Null
$d.Content.Find.Execute("446÷7=63, 5", $true, $false, $false, $false, $false, $true, 1, $false, "802÷8=100, 2", 2) | Out-Null
$d.Content.Find.Execute("460÷2=230, 0", $true, $false, $false, $false, $false, $true, 1, $false, "145÷8=18, 1", 2) | Out-Null
$d.Content.Find.Execute("810÷9=90, 0", $true, $false, $false, $false, $false, $true, 1, $false, "304÷2=152, 0", 2) | Out-Null
$d.Content.Find.Execute("354÷5=70, 4", $true, $false, $false, $false, $false, $true, 1, $false, "493÷4=123, 1", 2) | Out-Null
$d.Content.Find.Execute("527÷7=75, 2", $true, $false, $false, $false, $false, $true, 1, $false, "376÷9=41, 7", 2) | Out-Null
